$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 182-183; everything from the old row 182 onward
# shifts down by two rows (old row 182 -> new row 184, ... old row 287 ->
# new row 289). The sheet's used range grows from A1:R287 to A1:R289.
$ws.Rows("182:183").Insert()

# New row 182: weekly price record for "Región de O'Higgins" ($/caja 15 kilos)
$ws.Range("A182").Value = 6
$ws.Range("B182").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C182").Value = "Metropolitana"
$ws.Range("D182").Value = 44523
$ws.Range("E182").Value = 13
$ws.Range("F182").Value = 100112030
$ws.Range("G182").Value = "Poroto granado"
$ws.Range("H182").Value = "Sin especificar"
$ws.Range("I182").Value = "Primera"
$ws.Range("J182").Value = 95
$ws.Range("K182").Value = 35000
$ws.Range("L182").Value = 37000
$ws.Range("M182").Value = 35947
$ws.Range("N182").Value = "$/caja 15 kilos"
$ws.Range("O182").Value = "Región de O'Higgins"
$ws.Range("P182").Value = 2396
$ws.Range("Q182").Value = 15
$ws.Range("R182").Value = "Hortaliza"

# New row 183: weekly price record for "Perú" ($/malla 25 kilos)
$ws.Range("A183").Value = 6
$ws.Range("B183").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C183").Value = "Metropolitana"
$ws.Range("D183").Value = 44523
$ws.Range("E183").Value = 13
$ws.Range("F183").Value = 100112030
$ws.Range("G183").Value = "Poroto granado"
$ws.Range("H183").Value = "Sin especificar"
$ws.Range("I183").Value = "Primera"
$ws.Range("J183").Value = 250
$ws.Range("K183").Value = 45000
$ws.Range("L183").Value = 50000
$ws.Range("M183").Value = 48000
$ws.Range("N183").Value = "$/malla 25 kilos"
$ws.Range("O183").Value = "Perú"
$ws.Range("P183").Value = 1920
$ws.Range("Q183").Value = 25
$ws.Range("R183").Value = "Hortaliza"
